$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$xlPasteFormats = -4122
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$data = @(
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(1, 3),
    @(4, 5),
    @(6, 7),
    @(8, 8),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(4, 5),
    @(4, 5),
    @(5, 6),
    @(10, 10),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
